$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Introduction ")
$ws.Range("A1").Value = "test"
